$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 213, shifting existing rows 213:319 down to 214:320.
$ws.Rows("213:213").Insert()

# Populate the newly inserted row 213 with the new reading.
$ws.Range("A213").Value = 4
$ws.Range("B213").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C213").Value = "Los Lagos"
$ws.Range("D213").Value = 44523
$ws.Range("D213").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E213").Value = 10
$ws.Range("F213").Value = "Fruta"
$ws.Range("G213").Value = 100103
$ws.Range("H213").Value = "Frutos de hueso (carozo)"
$ws.Range("I213").Value = 100103006
$ws.Range("J213").Value = "Nectarín"
$ws.Range("K213").Value = "Early Glo"
$ws.Range("L213").Value = "Primera"
$ws.Range("M213").Value = 600
$ws.Range("N213").Value = 26000
$ws.Range("O213").Value = 27000
$ws.Range("P213").Value = 26500
$ws.Range("Q213").Value = "$/caja 14 kilos empedrada"
$ws.Range("R213").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S213").Value = 1893
$ws.Range("T213").Value = 14
